$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C5").Value = 79419.70291271698
$ws.Range("C6").Value = 77919.70291271698
$ws.Range("C7").Value = 67506.74747580943
$ws.Range("C9").Value = 16191.920719640257
$ws.Range("C10").Value = 63227.78219307674
$ws.Range("C11").Value = 61727.78219307674
$ws.Range("C12").Value = 19500.0
$ws.Range("C13").Value = 18000.0
$ws.Range("C14").Value = 43727.78219307675
$ws.Range("C15").Value = 42642.07325047115
$ws.Range("C16").Value = 397.078006105594
$ws.Range("C18").Value = 2919.1499999999996
$ws.Range("C19").Value = 40120.001256576754
$ws.Range("C20").Value = 23221.60255773338

$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C2").Value = 7941.5601221118795
$ws.Range("C3").Value = 9055.5
$ws.Range("D3").Value = 14.026713400891452
$ws.Range("C5").Value = 9055.499999999998
$ws.Range("C8").Value = 16818.0
$ws.Range("D8").Value = 111.77199116296094
$ws.Range("D9").Value = -6.6052527972599995
$ws.Range("D10").Value = 9.663590857309245
$ws.Range("D11").Value = 208.83100578325488
$ws.Range("D12").Value = 44.190811678384215
$ws.Range("C13").Value = 12152.0
$ws.Range("D13").Value = 53.01779263957078
$ws.Range("C14").Value = 7141.0
$ws.Range("D14").Value = -10.080640451022468
$ws.Range("C15").Value = 7463.0
$ws.Range("D15").Value = -6.026021521632921

$ws = $wb.Worksheets.Item("WING")
$ws.Range("C2").Value = 8418.05372943859
$ws.Range("C3").Value = 7831.5
$ws.Range("D3").Value = -6.967806909896137
$ws.Range("C5").Value = 7831.499999999999
$ws.Range("C8").Value = 7368.0
$ws.Range("D8").Value = -12.473830212872967
$ws.Range("C9").Value = 6522.0
$ws.Range("D9").Value = -22.523659154228756
$ws.Range("C10").Value = 9507.0
$ws.Range("D10").Value = 12.93584367076774
$ws.Range("C11").Value = 7929.0
$ws.Range("D11").Value = -5.8095819432505085

$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C2").Value = 913.2794140428662
$ws.Range("C3").Value = 866.0
$ws.Range("D3").Value = -5.17688380093576
$ws.Range("C5").Value = 865.9999999999999
$ws.Range("D8").Value = 99.50082876985572
$ws.Range("D9").Value = 112.42129956834253
$ws.Range("D10").Value = 30.080672106799444
$ws.Range("D11").Value = -188.47237631506226
$ws.Range("D12").Value = -12.294092291627647
$ws.Range("C13").Value = 513.0
$ws.Range("D13").Value = -43.82880068115479
$ws.Range("C14").Value = 626.0
$ws.Range("D14").Value = -31.45580745887504
$ws.Range("D15").Value = -7.366794105764034

$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C2").Value = 913.2794140428662
$ws.Range("C3").Value = 692.2857142857142
$ws.Range("D3").Value = -24.19781902001561
$ws.Range("C5").Value = 692.2857142857142
$ws.Range("D8").Value = 35.555447868870125
$ws.Range("D9").Value = 112.42129956834253
$ws.Range("D10").Value = -29.813374730253834
$ws.Range("D11").Value = -45.47123340977599
$ws.Range("D12").Value = -111.93501116131408
$ws.Range("C13").Value = 176.0
$ws.Range("D13").Value = -80.72878931751119
$ws.Range("C14").Value = 462.0
$ws.Range("D14").Value = -49.41307195846689

$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C2").Value = 1508.896423201257
$ws.Range("D3").Value = 7.098139749812221
$ws.Range("D10").Value = 5.772667723205633
$ws.Range("D11").Value = 6.965592547151562
$ws.Range("D12").Value = 8.556158979079466
$ws.Range("D17").Value = 5.772667723205633
$ws.Range("D18").Value = 6.965592547151562
$ws.Range("D19").Value = 8.556158979079466

$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("C2").Value = 6591.49490135286
$ws.Range("D3").Value = 21.874225097019256
$ws.Range("D11").Value = 30.046372306844297
$ws.Range("D12").Value = 17.909520015024952
$ws.Range("D13").Value = 17.666782969188564
$ws.Range("D18").Value = 30.046372306844297
$ws.Range("D19").Value = 17.909520015024952
$ws.Range("D20").Value = 17.666782969188564

$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C2").Value = 3256.039650065871
$ws.Range("C3").Value = 3160.31684344768
$ws.Range("D3").Value = -2.939853837967055
$ws.Range("C5").Value = 3160.3168434476797
$ws.Range("C9").Value = 3160.31684344768
$ws.Range("D9").Value = -2.9398538379670556
$ws.Range("C11").Value = 471.20536364093084
$ws.Range("C13").Value = 2689.1114798067492

$ws = $wb.Worksheets.Item("SYSTEMS")
$ws.Range("C2").Value = 10800.521766072157
$ws.Range("C3").Value = 8865.065365510069
$ws.Range("D3").Value = -17.92002685131349
$ws.Range("C4").Value = 8865.065365510069
$ws.Range("C8").Value = 8865.065365510069
$ws.Range("D8").Value = -17.92002685131349
$ws.Range("C11").Value = 391.6642079117032
$ws.Range("C13").Value = 391.66420791170316
$ws.Range("C21").Value = 1123.8004717186523
$ws.Range("C23").Value = 1123.800471718652
$ws.Range("C26").Value = 601.7732710517423
$ws.Range("C28").Value = 601.7732710517422
$ws.Range("C36").Value = 908.4284726885076
$ws.Range("C38").Value = 908.4284726885074
$ws.Range("C41").Value = 3736.658804044744
$ws.Range("C43").Value = 3736.6588040447436
